# PopUp ui changes. User result ready courses added
# Adds 3 new "ready" result rows (33, 34, 35) to the
# "EK5-PÇ Karşılama Yüzdeleri" worksheet, mirroring the existing
# row 32 pattern (same course/result template row, repeated for
# every ready course) with a couple of trailing cells trimmed off
# on the newly added rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$val1 = 0.9248000383377075
$valJ = 0.9371429085731506
$val2 = 0.6489999890327454
$val3 = 0.6800000071525574
$val4 = 0.6840000152587891
$val5 = 0.7344500422477722
$valAW = 0.7512667179107666
$valBB = 0.7553809285163879

# exact set of columns populated per new row (mirrors the diff exactly)
$rowsSpec = @{
  33 = @{
    g1 = @("E","F","G","H","I","J","K","L")
    g2 = @("P","Q","R","S","T","U","V")
    g3 = @("AA","AB","AC","AD","AE","AF","AG")
    g4 = @("AJ","AM","AN","AO","AP","AR","AU")
    g5 = @("AW","AX","AY","AZ","BA","BB","BC")
  }
  34 = @{
    g1 = @("E","F","G","H","I","J","K")
    g2 = @("P","Q","R","S","T","U","V")
    g3 = @("AA","AB","AC","AD","AE","AF","AG")
    g4 = @("AJ","AM","AN","AO","AP","AR","AU")
    g5 = @("AW","AX","AY","AZ","BA","BB","BC")
  }
  35 = @{
    g1 = @("E","F","G","H","I","J","K")
    g2 = @("P","Q","R","S","T","U","V")
    g3 = @("AA","AB","AC","AD","AE","AF","AG")
    g4 = @("AJ","AM","AN","AO","AP","AR","AU")
    g5 = @("AW","AX","AY","AZ","BA","BB","BC")
  }
}

foreach ($r in 33,34,35) {
  $spec = $rowsSpec[$r]

  $ws.Range("A$r").Value = "test"
  $ws.Range("B$r").Value = "test2"

  foreach ($col in $spec.g1) {
    $v = $val1
    if ($col -eq "J") { $v = $valJ }
    $ws.Range("$col$r").Value = $v
  }

  foreach ($col in $spec.g2) { $ws.Range("$col$r").Value = $val2 }

  foreach ($col in $spec.g3) { $ws.Range("$col$r").Value = $val3 }

  foreach ($col in $spec.g4) { $ws.Range("$col$r").Value = $val4 }

  foreach ($col in $spec.g5) {
    $v = $val5
    if ($col -eq "AW") { $v = $valAW }
    if ($col -eq "BB") { $v = $valBB }
    $ws.Range("$col$r").Value = $v
  }
}
